$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped crypto stats: the "Price" (D) and "Volume(1h)" (E)
# columns hold plain-text values (prices use "." as a thousands separator,
# e.g. "37.172.59", and volumes are padded strings like "  +1.84%  "), so
# every cell below is written back as text. "D" cells are additionally
# pinned to Text number format while the value is assigned, otherwise Excel
# would reinterpret them as numbers (dropping the thousands dots / trailing
# zeros); the style is reset straight back to "Normal" afterwards so no
# extra formatting lingers on the cell.

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '37.172.59'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.84%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.018.42'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +3.30%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '246.57'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '

# Row 6
$ws.Range("E6").Value = '  +1.14%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '60.14'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.390'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.30%  '

# Row 10
$ws.Range("E10").Value = '  +2.45%  '

# Row 11
$ws.Range("E11").Value = '  +1.05%  '

# Row 12
$ws.Range("E12").Value = '  +5.76%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '2.317.91'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +3.39%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.847'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.93%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '21.89'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.58%  '

# Row 16
$ws.Range("E16").Value = '  +3.28%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.022.20'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +3.20%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '37.160.58'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.97%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '70.27'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.52%  '

# Row 20
$ws.Range("E20").Value = '  +1.00%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.21'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.66%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '230.25'
$c.Style = "Normal"

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.55'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +4.34%  '

# Row 25
$ws.Range("E25").Value = '  -0.79%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.37'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +2.14%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '163.28'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.78%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.97%  '

# Row 29
$ws.Range("E29").Value = '  +2.42%  '

# Row 30
$ws.Range("E30").Value = '  +6.26%  '

# Row 31
$ws.Range("E31").Value = '  +0.75%  '

# Row 32
$ws.Range("E32").Value = '  +9.57%  '

# Row 33
$ws.Range("E33").Value = '  -0.05%  '

# Row 34
$ws.Range("E34").Value = '  +10.56%  '

# Row 35
$ws.Range("E35").Value = '  -0.29%  '

# Row 36
$ws.Range("E36").Value = '  +6.04%  '

# Row 37
$ws.Range("E37").Value = '  -0.01%  '

# Row 38
$ws.Range("E38").Value = '  +1.74%  '

# Row 39
$ws.Range("E39").Value = '  -1.99%  '

# Row 40
$ws.Range("E40").Value = '  +3.33%  '

# Row 41
$ws.Range("E41").Value = '  +0.65%  '

# Row 42
$ws.Range("E42").Value = '  +2.63%  '

# Row 43
$ws.Range("E43").Value = '  +1.55%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '16.67'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +5.23%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '91.25'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.85%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.376.88'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.21%  '

# Row 47
$ws.Range("E47").Value = '  +2.76%  '

# Row 48
$ws.Range("E48").Value = '  +3.57%  '

# Row 49
$ws.Range("E49").Value = '  +13.80%  '

# Row 50
$ws.Range("E50").Value = '  +1.56%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '46.13'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
